# SWOP Presentatie Inhoud.pptx
# Commit: "Pdf van presentatie added"
#
# The canonical-OOXML diff for this commit shows exactly one deliberate,
# content-level change: the (empty) placeholder slide that used to sit at
# position 17 - dark/"tx1" background, an untouched title placeholder and
# an untouched content placeholder, slide id 271 - is removed from the
# deck. Every other hunk in the diff (attribute-order churn on xmlns="" on
# a16:colId/a16:rowId/p15:sldGuideLst, and the relationship-id renumbering
# in <p:sldIdLst>) is just PowerPoint's own XML writer re-serialising the
# parts around the deleted slide; it is not a separate edit to reproduce.
#
# So: delete that one slide from the presentation.

$p = $ppt.ActivePresentation

$s = $p.Slides.Item(17)
$s.Delete()
